# Apply odds updates to Sheet1 as described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "V2"  = 1.57
    "G3"  = 1.73
    "H3"  = 3.7
    "I3"  = 5
    "J3"  = 2.38
    "K3"  = 2.1
    "Q3"  = 2.2
    "R3"  = 1.65
    "V3"  = 1.67
    "W3"  = 5.5
    "X3"  = 7
    "AA3" = 15
    "AC3" = 8
    "AH3" = 11
    "AJ3" = 17
    "AN3" = 3.5
    "M4"  = 1.1
    "N4"  = 7
    "V4"  = 1.67
    "U5"  = 2.63
    "V5"  = 1.44
    "M12" = 1.1
    "N12" = 7
    "Q12" = 2.4
    "R12" = 1.53
    "U12" = 2.38
    "V12" = 1.53
    "G13" = 3.75
    "H13" = 3.3
    "I13" = 1.95
    "J13" = 4.33
    "L13" = 2.63
    "M13" = 1.05
    "N13" = 11
    "U13" = 1.8
    "V13" = 1.91
    "AB13" = 41
    "AE13" = 15
    "AF13" = 51
    "AH13" = 7.5
    "AI13" = 9.5
    "AK13" = 17
    "AL13" = 15
    "AO13" = 21
    "AP13" = 29
    "AZ13" = 34
    "G14" = 2.35
    "H14" = 3.15
    "I14" = 2.77
    "J14" = 3
    "L14" = 3.4
    "O14" = 1.36
    "S14" = 1.42
    "T14" = 2.65
    "U14" = 1.83
    "V14" = 1.87
    "W14" = 7.4
    "X14" = 11.25
    "Y14" = 9.5
    "AA14" = 21
    "AE14" = 15
    "AF14" = 75
    "AJ14" = 10.5
    "AK14" = 32
    "AN14" = 4.25
    "AP14" = 21
    "AQ14" = 55
    "AT14" = 2.65
    "AU14" = 7.2
    "AY14" = 24
    "BA14" = 110
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
